$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Marseille: updated stats after a new match (loss)
$ws.Range("C5").Value = 23
$ws.Range("F5").Value = 7
$ws.Range("H5").Value = 31
$ws.Range("I5").Value = 17
$ws.Range("K5").Value = 1.74
$ws.Range("L5").Value = "W D L D L"

# Rows 12/13 - Brest and Angers swap order, with Brest's stats updated
$ws.Range("B12").Value = "Brest"
$ws.Range("C12").Value = 23
$ws.Range("E12").Value = 6
$ws.Range("G12").Value = 31
$ws.Range("H12").Value = 34
$ws.Range("I12").Value = -3
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 1.3
$ws.Range("L12").Value = "L D W D W"
$ws.Range("M12").Value = 14461
$ws.Range("N12").Value = "Romain Del Castillo - 7"
$ws.Range("O12").Value = "Grégoire Coudert"

$ws.Range("B13").Value = "Angers"
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = 5
$ws.Range("G13").Value = 22
$ws.Range("H13").Value = 27
$ws.Range("J13").Value = 29
$ws.Range("K13").Value = 1.32
$ws.Range("L13").Value = "L D W W L"
$ws.Range("M13").Value = 12232
$ws.Range("N13").Value = "Sidiki Cherif - 4"
$ws.Range("O13").Value = "Hervé Koffi"
